$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.79%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'29.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'9.38%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.173"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'2.08%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05706"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.32%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.601"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.96%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'3.067"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.20%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.8590"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'4.67%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D9").Value = "'0.8735"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'4.24%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1368"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.26%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07097"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.85%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.02877"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.75%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09397"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.00%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001511"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.23%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "'0.04176"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'2.13%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006085"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.20%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'3,766.80%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'-0.72%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.281"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.31%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "One"
$ws.Range("C20").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D20").Value = "'0.01028"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1,607.92%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'0.72%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.03296"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'3.70%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'0.35%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'3.472"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-3.12%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.1381"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.52%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.005036"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'27.48%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'0.12%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001211"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'23.57%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.03752"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.76%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.005773"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.95%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1072"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.76%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002101"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-10.02%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01022"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'9.01%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005178"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-0.71%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.08%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.07004"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-30.97%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002569"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-0.89%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.08%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.08%"
$ws.Range("E50").Style = "Normal"
